# "actualizacion: imagenes al catalogo"
# Appends 7 new product codes to column A of the active sheet (MEJORAR),
# right after the existing data (which currently ends at row 44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "SAHANAF1H",
    "SAHANAF2H",
    "EVOL0043",
    "EVOL3420",
    "EVOL5530",
    "EVOL4755",
    "EVOL4753"
)

$startRow = $ws.Cells(1, 1).End(-4121).Row + 1   # xlDown = -4121; last used row in col A, +1

$row = $startRow
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row++
}

# Match the author's final selection (the next empty cell after the new rows).
$null = $ws.Cells.Item($row, 1).Select()
